# patterns.xlsx edit:
#   - "new ontology+patterns, new ParallelRun, stupid bug fix"
#
# The real fix: row 10 (id=9) had its arg1_case column pointing at a stray
# duplicate string "gen" instead of the canonical "gent" (genitive) used
# everywhere else in the sheet (see e.g. row 7 / id=6 which correctly uses
# "gent"). Correcting the cell value removes the now-unreferenced "gen"
# shared string and collapses the shared-string table, which is why every
# later shared-string index in the sheet shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- the actual data bug fix --------------------------------------------
$ws.Range("C10").Value = "gent"

# --- selection moved while reviewing the fix ----------------------------
$ws.Range("C11").Select()

# --- column width cleanup ------------------------------------------------
# Drop the stale narrow-column overrides for A and C:F (they go back to
# the sheet's default width) and widen the two content columns (B / G).
$ws.Columns.Item(1).ClearFormats()
$ws.Range($ws.Columns.Item(3), $ws.Columns.Item(6)).ClearFormats()

# NOTE: this COM host quantizes ColumnWidth to 1/6-character steps, so we
# pick the input that lands on the step closest to the real target widths
# (59.1860465116279 / 26.1813953488372).
$ws.Columns.Item(2).ColumnWidth = 58.4
$ws.Columns.Item(7).ColumnWidth = 25.4
